$d = $word.ActiveDocument
$shape = $d.InlineShapes.Item(1)
$picStart = $shape.Range.Start
$shape.Delete()

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/Retail_Unit_Size.png?h=60%25&w=60%25"
$target = $d.Range($picStart, $picStart)
$d.Hyperlinks.Add($target, $url, $null, $null, $url)
Write-Host "OK"
